$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 47 (continuation of the work log for 4.4.2020) ---

# A47: date label "4.4.2020" -- already exists as a shared string, keep it
#      stored as text (not auto-converted to a date serial number).
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = "4.4.2020"

# B47 / C47: from / to time-of-day values
$ws.Range("B47").Value = 0.61458333333333337
$ws.Range("C47").Value = 0.64236111111111105

# D47: elapsed time formula (C47-B47)
$ws.Range("D47").FormulaR1C1 = "=RC[-1]-RC[-2]"

# E47: task -- "IO Debounce" (existing shared string)
$ws.Range("E47").Value = "IO Debounce"

# F47: note -- "Improve Testbench" (new shared string)
$ws.Range("F47").Value = "Improve Testbench"

# G47: project -- "Use added packages" (new shared string)
$ws.Range("G47").Value = "Use added packages"

# Re-apply the formatting of the row directly above (row 46) onto the new
# row so number formats / alignment match the rest of the table, without
# touching the values/formulas we just entered.
$ws.Range("A46:G46").Copy()
$ws.Range("A47:G47").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Add an extra (otherwise empty) formatted cell at D52 ---
$ws.Range("D52").NumberFormat = "h:mm"

# --- Update the selection so it highlights G47, the last-edited cell ---
$ws.Range("G47").Select() | Out-Null
